# Apply "resolve and classify+summarise" re-run results after mapping file changes.
$wb = $excel.ActiveWorkbook

# --- Sheet "Range Status": zero out species counts, clear percentage column ---
$wsRange = $wb.Worksheets.Item("Range Status")
$wsRange.Range("B2:B7").Value = 0
$wsRange.Range("C2:C7").ClearContents()

# --- Sheet "Species qualification": Range Analysis species count -> 0 ---
$wsQual = $wb.Worksheets.Item("Species qualification")
$wsQual.Range("B5").Value = 0

# --- Sheet "High Priority break-up": remove Range row, turn old IUCN row into row 3 ---
$wsBreak = $wb.Worksheets.Item("High Priority break-up")
$wsBreak.Range("E2").Value = 2.6

$wsBreak.Range("A3").Value = "IUCN"
$wsBreak.Range("B3").Value = 38
$wsBreak.Range("C3").Value = 97.40000000000001
$wsBreak.Range("D3").Value = 38
$wsBreak.Range("E3").Value = 97.40000000000001

# delete old row 4 (previous IUCN row), shifting rows up
$wsBreak.Rows.Item(4).Delete()
